$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Studies Searches")
$vals = 18.0,18.140625,19.0,20.0,25.0,27.0,27.5,27.5703125,27.57,28.0,29.0,11.0,11.28515625,11.3,15.0
$c = 1
foreach ($v in $vals) {
  $ws.Columns.Item($c).ColumnWidth = $v
  $c = $c + 1
}
